$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# clear rows 3:7 first
$ws.Range("A3:C7").Clear()

$ws.Range("A1").Value = "№"
$ws.Range("B1").Value = "Посыл / Заповедь"
$ws.Range("C1").Value = "Время (мск)"

$ws.Range("B2").Value = "U2FsdGVkX1+uSMZWPh9qT+dXpFUxouwT+4p0znZbKwHwZnJH37GSiCeCmdRsRdYBpFjz6yycSGEHbXPiOZQSOgx25Ak70aJLQqHrB98ER2HZyEE1omZZ3yIMIgqEqcXF3s1kW1Fwo7RiKjFJ0kzFfVergxEiftrNFpSNdTib9ol+4Yt5vs8BGAuVZMfZhWq6aK3dyvHQk+ZJLCt0jDT/JZ64obTFWcIKQ2p6mc0gTlji+Rmh7bnonmKScPBEL4OjfMJId+hWCfWJ6Yjp2ZqsKVKWMSrgDfSxFv3FP50XMyjL1FP7S2rQTc9+JjePWzSi7RHlbX1+mRuCaqmY164GROmMnTqSaSkS8odllccGmDoc5H6WhVfKE6TV3r6cuZWX7u6sdPEvVf3w2IsoQi7ZB0AebPrcf88T8KbjFeFcGsmSgQVpQ0+u8cbJ7WlPMHQxBKh7kIxf+pGnbmiVLBjnfwOo/cSlyxxNEHSveN3/Xu5C6EJCe9+yiwtKhKFynSjClM/56la5ZAC+mNeEfqmZtMAbXC3d4PWE9QcJ9Oqa7JdnTby8HC5YYEmCThqLBd7Lcw4Lk/JJRUPg9MM8kJ+3W0RM4H33vp4LQn0MrZKC2mDmmXUtxNpHWJTuErBOjRoydThNbnoIJ8nMsm+dBHpv5oj9JeDd77mw524WEs8j0p1TUOWghTGfj76BsZEGJ6P1PFSSNIB9y9+3rkzCebzEGbhH2+eL8yomOWJxSAxJqGA4EPTqTbw1Dl3uH58VajHho1d3buD0fS297IPu6pLJKuTaUlvlgdLxtdnlhQkW+LVZKyOzyIh3B8MKlkdHOQNzY4ib+sn0xQJj4mNbVtZ5slwRJko3Iklrh666byXqUniklEYABJh80MJ0ojjyWld1F6GJRVxX6ZwndnFoOGMDm0JwpMejwvqnoDMhR/G5wgqDkG/JnbaQHnPj4W7AiWLkln+u1fADEi6Vyinrxg3x1CZWfCSorW3Jg+nQN1H4K31XgDbhqNUcr28qxOW//hzyQvlWN4DY6NngfuUdIUWx2hmp5sx0Ev93IwRGEJ2UxK5due/vdDp7rOeUtxdLhP9vE5pwYqgISkt7N/9aUt3r7w887VdlBuDWg2YgG4+0cLVnoAJhptrKhXljErWA4gLrUTdnbotHwulTY3WBHP0YKKGLAVVa2hA8vB9U9fmm+lgyukO5bpHHeKrbUCjupZgA8MtcO4JkLePicUCh0Dd2xx0F8VRXC5ZOnF6h9S8myqnbxqvMqyYgGJ3agqFOfjYo+N3NETW9tZQUaOfWMeSFGX8BqVhhzWwT0fDN4OsealA="
$ws.Range("C2").Value = "U2FsdGVkX1/gIUGEizncYlSscEtCg92aHZMTBJALFwc="

$ws.Range("B11").Select()
